$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "Autotest_1"
$ws.Range("P1").Value = "Appellant_Valid"
$ws.Range("Q1").Value = "Appellant_Invalid"
$ws.Range("Q2").Value = "Autotest_5"

$ws.Range("A2").Select()
